$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '{"username":"cisco","password":"cisco","timeout":10,"type":"cisco"}'
$ws.Range("D6").Value = '{"username":"cisco","password":"cisco","timeout":10,"type":"cisco"}'

$ws.Range("D6").Select()
